$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = -5
$ws.Range("F11").Value = -1
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = 6
$ws.Range("F16").Value = -4
$ws.Range("F17").Value = -2
$ws.Range("F21").Value = -6
$ws.Range("F22").Value = -8
$ws.Range("F24").Value = -7

$wb.Save()
